$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1596
$ws.Range("G2").Value = 70
$ws.Range("F3").Value = 3330
$ws.Range("F5").Value = 759
$ws.Range("F6").Value = 2373
$ws.Range("F8").Value = 428
$ws.Range("F9").Value = 257
$ws.Range("F10").Value = 154
$ws.Range("F11").Value = 376
$ws.Range("F14").Value = 231
$ws.Range("F16").Value = 281
$ws.Range("F17").Value = 4924
$ws.Range("F18").Value = 30
$ws.Range("F19").Value = 1390
$ws.Range("F20").Value = 3623
$ws.Range("F21").Value = 167
$ws.Range("F22").Value = 214
$ws.Range("F23").Value = 3933
$ws.Range("F24").Value = 5300
$ws.Range("F27").Value = 581
$ws.Range("F28").Value = 3392
$ws.Range("F29").Value = 396
$ws.Range("F31").Value = 154
$ws.Range("F32").Value = 98
$ws.Range("F34").Value = 1227
$ws.Range("F35").Value = 43
$ws.Range("F36").Value = 64
$ws.Range("F37").Value = 1455
$ws.Range("F38").Value = 151
$ws.Range("F39").Value = 1437
$ws.Range("F41").Value = 935
$ws.Range("F42").Value = 922
$ws.Range("F43").Value = 532
$ws.Range("F45").Value = 2466
$ws.Range("F46").Value = 90
$ws.Range("F47").Value = 189
$ws.Range("F48").Value = 378
$ws.Range("F49").Value = 3770

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 1036

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2627

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2627
$ws.Range("F3").Value = 1596
$ws.Range("G3").Value = 70
$ws.Range("F4").Value = 3330
$ws.Range("F6").Value = 759
$ws.Range("F8").Value = 2373
$ws.Range("F10").Value = 428
$ws.Range("F11").Value = 257
$ws.Range("F12").Value = 1036
$ws.Range("F13").Value = 154
$ws.Range("F14").Value = 376
$ws.Range("F17").Value = 231
$ws.Range("F19").Value = 281
$ws.Range("F20").Value = 4924
$ws.Range("F22").Value = 1390
$ws.Range("F23").Value = 3933
$ws.Range("F24").Value = 5300
$ws.Range("F27").Value = 581
$ws.Range("F28").Value = 3392
$ws.Range("F29").Value = 396
$ws.Range("F31").Value = 154
$ws.Range("F32").Value = 98
$ws.Range("F33").Value = 1227
$ws.Range("F34").Value = 43
$ws.Range("F35").Value = 64
$ws.Range("F36").Value = 1455
$ws.Range("F37").Value = 1437
$ws.Range("F38").Value = 935
$ws.Range("F39").Value = 532
$ws.Range("F43").Value = 2466
$ws.Range("F45").Value = 90
$ws.Range("F46").Value = 189
$ws.Range("F47").Value = 378
$ws.Range("F49").Value = 3770
